$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column keeps its text formatting so numeric-looking
# values (e.g. "306.10") are not coerced into floating point numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "42.358.83"
$ws.Range("E2").Value = "  +1.24%  "
$ws.Range("D3").Value = "2.269.39"
$ws.Range("E3").Value = "  -0.15%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "306.10"
$ws.Range("E5").Value = "  +0.63%  "
$ws.Range("D6").Value = "97.56"
$ws.Range("E6").Value = "  +5.19%  "
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "0.491"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").Value = "35.68"
$ws.Range("E10").Value = "  +9.51%  "
$ws.Range("E11").Value = "  -0.12%  "
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "6.65"
$ws.Range("E13").Value = "  -0.68%  "
$ws.Range("D14").Value = "2.624.24"
$ws.Range("E14").Value = "  +0.01%  "
$ws.Range("D15").Value = "14.36"
$ws.Range("E15").Value = "  +0.37%  "
$ws.Range("D16").Value = "2.278.10"
$ws.Range("E16").Value = "  -0.40%  "
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "42.258.56"
$ws.Range("E18").Value = "  +1.18%  "
$ws.Range("D19").Value = "12.51"
$ws.Range("E19").Value = "  -0.96%  "
$ws.Range("D20").Value = "0.0₃0910"
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  +0.20%  "
$ws.Range("D22").Value = "67.58"
$ws.Range("E22").Value = "  +0.60%  "
$ws.Range("D23").Value = "238.43"
$ws.Range("E23").Value = "  -2.05%  "
$ws.Range("D24").Value = "2.59"
$ws.Range("E24").Value = "  +0.30%  "
$ws.Range("E25").Value = "  +1.67%  "
$ws.Range("E26").Value = "  -0.58%  "
$ws.Range("D27").Value = "23.79"
$ws.Range("E27").Value = "  -1.02%  "
$ws.Range("D28").Value = "37.46"
$ws.Range("E28").Value = "  +6.19%  "
$ws.Range("D29").Value = "9.52"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "2.10"
$ws.Range("E30").Value = "  +1.53%  "
$ws.Range("D31").Value = "160.18"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("D32").Value = "5.25"
$ws.Range("E32").Value = "  +0.06%  "
$ws.Range("E33").Value = "  +0.06%  "
$ws.Range("D34").Value = "3.15"
$ws.Range("E34").Value = "  +4.51%  "
$ws.Range("D36").Value = "17.11"
$ws.Range("E36").Value = "  +0.25%  "
$ws.Range("E37").Value = "  -0.66%  "
$ws.Range("D38").Value = "2.35"
$ws.Range("E38").Value = "  -0.55%  "
$ws.Range("E39").Value = "  +1.69%  "
$ws.Range("D40").Value = "0.114"
$ws.Range("E40").Value = "  -1.49%  "
$ws.Range("E41").Value = "  +3.82%  "
$ws.Range("E42").Value = "  +14.38%  "
$ws.Range("D43").Value = "1.994.00"
$ws.Range("E43").Value = "  -0.69%  "
$ws.Range("D44").Value = "0.0287"
$ws.Range("E44").Value = "  +2.16%  "
$ws.Range("D45").Value = "19.00"
$ws.Range("E45").Value = "  -2.71%  "
$ws.Range("B46").Value = "FraxShare"
$ws.Range("C46").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D46").Value = "9.99"
$ws.Range("E46").Value = "  -3.03%  "
$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "2.94"
$ws.Range("E47").Value = "  +1.14%  "
$ws.Range("D48").Value = "53.25"
$ws.Range("E48").Value = "  +0.68%  "
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("D50").Value = "72.06"
$ws.Range("E50").Value = "  -0.34%  "
$ws.Range("D51").Value = "91.42"
$ws.Range("E51").Value = "  +0.22%  "
